# Replace the stale AQL 7.x stack trace embedded in the bold red run with the
# updated AQL 8.x stack trace. The leading "    <---divOp(...)...failed:" line
# is left untouched; only the trace text starting at "divOp(" (the stack frames)
# is replaced, so the Find locates that span using a wildcard match anchored on
# stable start/end tokens.
$d = $word.ActiveDocument
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("divOp*RemoteTestRunner*210)", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text block"
}
# The matched range ends right before a literal "`n" that is baked into the
# w:t run (followed by a separate <w:br/>). Extend by one character so that
# assigning .Text below does not leave a duplicate trailing newline.
$rng.End = $rng.End + 1

$newText = @'
divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:
	java.lang.ArithmeticException: / by zero
		at org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)
		at java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)
		at java.base/java.lang.reflect.Method.invoke(Method.java:580)
		at org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:170)
		at org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:231)
		at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCallService(EvaluationServices.java:122)
		at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCall(EvaluationServices.java:237)
		at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:170)
		at org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:244)
		at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:135)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:241)
		at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:135)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:146)
		at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:53)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:2362)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1)
		at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:134)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:2349)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)
		at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:2373)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1)
		at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:134)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:2349)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)
		at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:350)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)
		at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:334)
		at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:878)
		at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:499)
		at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:391)
		at java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)
		at java.base/java.lang.reflect.Method.invoke(Method.java:580)
		at org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:59)
		at org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)
		at org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)
		at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)
		at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)
		at org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)
		at org.junit.runners.BlockJUnit4ClassRunner$1.evaluate(BlockJUnit4ClassRunner.java:100)
		at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)
		at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)
		at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)
		at org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)
		at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)
		at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)
		at org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)
		at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)
		at org.junit.runners.ParentRunner.run(ParentRunner.java:413)
		at org.junit.runners.Suite.runChild(Suite.java:128)
		at org.junit.runners.Suite.runChild(Suite.java:27)
		at org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)
		at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)
		at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)
		at org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)
		at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)
		at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)
		at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)
		at org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)
		at org.junit.runners.ParentRunner.run(ParentRunner.java:413)
		at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:93)
		at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:40)
		at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:529)
		at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:757)
		at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:452)
		at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:210)
	
'@

$rng.Text = $newText
